$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force text number-format so numeric-looking strings (e.g. "5.80", "0.999")
    # are not coerced to numbers / stripped of trailing zeros, then restore the
    # default "Normal" style so no stray style index is left on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "98.769.40"
Set-TextValue $ws.Range("E2") "  +0.54%  "

Set-TextValue $ws.Range("D3") "3.317.59"
Set-TextValue $ws.Range("E3") "  -1.56%  "

Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "255.74"
Set-TextValue $ws.Range("E5") "  +0.30%  "

Set-TextValue $ws.Range("D6") "626.57"
Set-TextValue $ws.Range("E6") "  +0.49%  "

Set-TextValue $ws.Range("D7") "1.48"
Set-TextValue $ws.Range("E7") "  +23.17%  "

Set-TextValue $ws.Range("D8") "0.418"
Set-TextValue $ws.Range("E8") "  +8.23%  "

Set-TextValue $ws.Range("B9") "USDC"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  +0.02%  "

Set-TextValue $ws.Range("B10") "Cardano"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D10") "1.02"
Set-TextValue $ws.Range("E10") "  +24.60%  "

Set-TextValue $ws.Range("D11") "3.314.73"

Set-TextValue $ws.Range("E12") "  +2.76%  "

Set-TextValue $ws.Range("D13") "41.29"
Set-TextValue $ws.Range("E13") "  +14.93%  "

Set-TextValue $ws.Range("D14") "98.506.68"
Set-TextValue $ws.Range("E14") "  +0.53%  "

Set-TextValue $ws.Range("E15") "  +2.83%  "

Set-TextValue $ws.Range("D16") "3.945.66"
Set-TextValue $ws.Range("E16") "  -1.13%  "

Set-TextValue $ws.Range("E17") "  -1.62%  "

Set-TextValue $ws.Range("D18") "3.298.51"
Set-TextValue $ws.Range("E18") "  -2.17%  "

Set-TextValue $ws.Range("E19") "  -5.09%  "

Set-TextValue $ws.Range("D20") "15.71"
Set-TextValue $ws.Range("E20") "  +4.87%  "

Set-TextValue $ws.Range("D21") "6.43"
Set-TextValue $ws.Range("E21") "  +8.69%  "

Set-TextValue $ws.Range("D22") "485.93"
Set-TextValue $ws.Range("E22") "  +0.76%  "

Set-TextValue $ws.Range("D23") "9.47"
Set-TextValue $ws.Range("E23") "  +2.51%  "

Set-TextValue $ws.Range("D24") "0.0000203"
Set-TextValue $ws.Range("E24") "  -2.10%  "

Set-TextValue $ws.Range("D25") "5.80"
Set-TextValue $ws.Range("E25") "  +0.05%  "

Set-TextValue $ws.Range("D26") "0.343"
Set-TextValue $ws.Range("E26") "  +36.01%  "

Set-TextValue $ws.Range("D27") "89.45"

Set-TextValue $ws.Range("D28") "12.18"
Set-TextValue $ws.Range("E28") "  +1.08%  "

Set-TextValue $ws.Range("D29") "3.492.42"
Set-TextValue $ws.Range("E29") "  -1.46%  "

Set-TextValue $ws.Range("D30") "0.151"
Set-TextValue $ws.Range("E30") "  +20.34%  "

Set-TextValue $ws.Range("E31") "  -0.11%  "

Set-TextValue $ws.Range("E32") "  +2.93%  "

Set-TextValue $ws.Range("D33") "10.62"
Set-TextValue $ws.Range("E33") "  +14.58%  "

Set-TextValue $ws.Range("E34") "  +0.21%  "

Set-TextValue $ws.Range("D35") "28.09"
Set-TextValue $ws.Range("E35") "  +2.53%  "

Set-TextValue $ws.Range("D36") "0.482"
Set-TextValue $ws.Range("E36") "  +7.48%  "

Set-TextValue $ws.Range("E37") "  -0.69%  "

Set-TextValue $ws.Range("D38") "7.37"
Set-TextValue $ws.Range("E38") "  +0.17%  "

Set-TextValue $ws.Range("E39") "  +0.78%  "

Set-TextValue $ws.Range("D40") "498.59"
Set-TextValue $ws.Range("E40") "  -5.39%  "

Set-TextValue $ws.Range("E41") "  -0.28%  "

Set-TextValue $ws.Range("D42") "3.90"
Set-TextValue $ws.Range("E42") "  -0.80%  "

Set-TextValue $ws.Range("E43") "  -1.28%  "

Set-TextValue $ws.Range("D44") "0.792"
Set-TextValue $ws.Range("E44") "  +0.07%  "

Set-TextValue $ws.Range("E45") "  +0.00%  "

Set-TextValue $ws.Range("E46") "  -2.09%  "

Set-TextValue $ws.Range("D47") "160.28"
Set-TextValue $ws.Range("E47") "  -0.57%  "

Set-TextValue $ws.Range("E48") "  +1.22%  "

Set-TextValue $ws.Range("D49") "0.860"
Set-TextValue $ws.Range("E49") "  +7.78%  "

Set-TextValue $ws.Range("D50") "4.77"
Set-TextValue $ws.Range("E50") "  +5.01%  "

Set-TextValue $ws.Range("D51") "7.32"
Set-TextValue $ws.Range("E51") "  +13.76%  "

